$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

# Copy the formatting (fill style) of row 16 down into the new rows 19:22
$ws.Range("A16:L16").Copy()
$ws.Range("A19:L22").PasteSpecial($xlPasteFormats)

$rows = @(
    @{ Row = 19; G = 0 },
    @{ Row = 20; G = 0.25 },
    @{ Row = 21; G = 0.5 },
    @{ Row = 22; G = 0.75 }
)

foreach ($row in $rows) {
    $r = $row.Row
    $ws.Range("A$r").Value = "zeroshot huang combined with own (w/o marketing)"
    $ws.Range("D$r").Value = 4000
    $ws.Range("E$r").Value = 200
    $ws.Range("F$r").Value = 5
    $ws.Range("G$r").Value = $row.G
    $ws.Range("H$r").Value = 4000
    $ws.Range("I$r").Value = "null"
    $ws.Range("J$r").Value = "yes"
}

$ws.Range("I26").Select()
